$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Year" header to A1 (B1 already contains "PpD")
$ws.Range("A1").Value = "Year"

# Update selection to E4 as recorded in the saved workbook view
$ws.Range("E4").Select()
